# Update the "Förändrad" (changed) date column C for rows 2-5
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45212
$ws.Range("C3").Value = 45212
$ws.Range("C4").Value = 45212
$ws.Range("C5").Value = 45212
